# Apply property-sheet header/column normalization to the
# insurance (保險, sheet 6) and business-investment (事業投資, sheet 7)
# worksheets, per commit "#5: insurance, claim, debt, investment done".
#
# Both sheets originally used raw data values as their header row
# (a bug) and were missing the common trailing metadata columns
# (property_category, category, date, legislator_name, legislator_id,
# source_file, index) that the rest of the workbook's sheets already
# have. This script fixes the headers and appends the missing columns
# with the correct values for every existing data row.

$wb = $excel.ActiveWorkbook

# A cell elsewhere in the workbook that already holds the literal text
# "2013-12-30" as a (non-date-formatted) shared string. Excel's COM
# layer auto-converts a literal ISO-date-looking string assigned via
# .Value into a real date serial number, so for that one value we
# instead copy an existing text cell to guarantee it lands back in
# the workbook as plain text, matching every other sheet.
$dateSourceCell = $wb.Worksheets.Item("土地").Cells.Item(2, 11)

# ---------------------------------------------------------------
# Sheet 6: 保險 (insurance)
# ---------------------------------------------------------------
$ws6 = $wb.Worksheets.Item("保險")

# Extend formatting (bold header style) from existing header cells
# into the new header columns F1:K1, then overwrite with real text.
$ws6.Range("B1:E1").Copy($ws6.Range("F1:K1"))

$ws6.Range("B1").Value = "company"
$ws6.Range("C1").Value = "name"
$ws6.Range("D1").Value = "owner"
$ws6.Range("E1").Value = "property_category"
$ws6.Range("F1").Value = "category"
$ws6.Range("G1").Value = "date"
$ws6.Range("H1").Value = "legislator_name"
$ws6.Range("I1").Value = "legislator_id"
$ws6.Range("J1").Value = "source_file"
$ws6.Range("K1").Value = "index"

# Data rows 2-7: fix B/C columns (previously duplicated header text
# into the data rows) and append F:K metadata columns.
$sheet6Rows = @(
    @{ Row = 2;  Index = 80; Company = "富邦人壽"; Name = "20LPL安泰分紅終身壽險"; Owner = "蘇震清" },
    @{ Row = 3;  Index = 81; Company = "富邦人壽"; Name = "20LPL安泰分紅終身壽險"; Owner = "廖靖汝" },
    @{ Row = 4;  Index = 82; Company = "富邦人壽"; Name = "20LPL安泰分紅終身壽險"; Owner = "蘇震清" },
    @{ Row = 5;  Index = 83; Company = "富邦人壽"; Name = "20LPL安泰分紅終身壽險"; Owner = "蘇震清" },
    @{ Row = 6;  Index = 84; Company = "富邦人壽"; Name = "安泰喬壽還本終身壽險"; Owner = "蘇震清" },
    @{ Row = 7;  Index = 85; Company = "富邦人壽"; Name = "安泰喬壽還本終身壽險"; Owner = "蘇震清" }
)

foreach ($r in $sheet6Rows) {
    $row = $r.Row
    $ws6.Cells.Item($row, 2).Value = $r.Company
    $ws6.Cells.Item($row, 3).Value = $r.Name
    $ws6.Cells.Item($row, 4).Value = $r.Owner
    $ws6.Cells.Item($row, 5).Value = "insurance"
    $ws6.Cells.Item($row, 6).Value = "normal"
    $dateSourceCell.Copy($ws6.Cells.Item($row, 7))
    $ws6.Cells.Item($row, 8).Value = "蘇震清"
    $ws6.Cells.Item($row, 9).Value = 1718
    $ws6.Cells.Item($row, 10).Value = "tmpb3b61"
    $ws6.Cells.Item($row, 11).Value = $r.Index
}

# ---------------------------------------------------------------
# Sheet 7: 事業投資 (business investment)
# ---------------------------------------------------------------
$ws7 = $wb.Worksheets.Item("事業投資")

# Extend formatting (bold header style) from existing header cells
# into the new header columns H1:N1, then overwrite with real text.
$ws7.Range("B1:G1").Copy($ws7.Range("H1:N1"))

$ws7.Range("B1").Value = "owner"
$ws7.Range("C1").Value = "company"
$ws7.Range("D1").Value = "address"
$ws7.Range("E1").Value = "total"
$ws7.Range("F1").Value = "register_date"
$ws7.Range("G1").Value = "register_reason"
$ws7.Range("H1").Value = "property_category"
$ws7.Range("I1").Value = "category"
$ws7.Range("J1").Value = "date"
$ws7.Range("K1").Value = "legislator_name"
$ws7.Range("L1").Value = "legislator_id"
$ws7.Range("M1").Value = "source_file"
$ws7.Range("N1").Value = "index"

# Data row 2: fix C/D columns (previously duplicated header text into
# the data row) and append H:N metadata columns.
$ws7.Cells.Item(2, 2).Value = "廖靖汝"
$ws7.Cells.Item(2, 3).Value = "南島休閒育樂股份有限公司"
$ws7.Cells.Item(2, 4).Value = "高雄市精富路148號"
$ws7.Cells.Item(2, 5).Value = 1000000
$ws7.Cells.Item(2, 6).Value = "95年08月23日"
$ws7.Cells.Item(2, 7).Value = "投資"
$ws7.Cells.Item(2, 8).Value = "investment"
$ws7.Cells.Item(2, 9).Value = "normal"
$dateSourceCell.Copy($ws7.Cells.Item(2, 10))
$ws7.Cells.Item(2, 11).Value = "蘇震清"
$ws7.Cells.Item(2, 12).Value = 1718
$ws7.Cells.Item(2, 13).Value = "tmpb3b61"
$ws7.Cells.Item(2, 14).Value = 98

Write-Output "sheet6/sheet7 headers and metadata columns updated"
